$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading-percent values per row (row index -> values for columns B,C,E,F,G,H,M,N,O)
# Rows 2-25 correspond to time steps 0-23 in column A
$newValues = @{
    2 = @(10.33992315719617, 8.777269625212423, 12.6172380270932, 16.86991607391245, 17.64171322286508, 11.35281987289404, 13.38866338775711, 16.00341270652173, 15.79316265678191)
    3 = @(9.73012301025404, 8.539781048769274, 12.49652207463946, 15.89584955866815, 17.58493340090801, 11.39521390737662, 13.11164444333345, 16.0406846585789, 15.84611711151641)
    4 = @(9.333670343419312, 8.389694039751355, 12.42677375572341, 15.26997757108489, 17.55933004182025, 11.42350618813135, 12.94134546367982, 16.0653452144461, 15.88335191161698)
    5 = @(9.166612735328403, 8.327523494905057, 12.39948182964083, 15.008197319934, 17.55122957547699, 11.43560354757475, 12.8719947106504, 16.07584168600638, 15.89970726665762)
    6 = @(9.138542181223027, 8.317141164748062, 12.39501921359238, 14.96433081551589, 17.5500254614865, 11.43764658824928, 12.86048489793031, 16.077611639792, 15.90249428741372)
    7 = @(9.331439553128869, 8.38885958235954, 12.42640106768547, 15.26647399323133, 17.55921134661254, 11.42366703852988, 12.94040984947035, 16.06548496239964, 15.88356770731445)
    8 = @(10.13425774965859, 8.696306533978508, 12.57472771025954, 16.5399640634477, 17.62021789740195, 11.36696733464614, 13.29324961704647, 16.01589597389943, 15.81043858496723)
    9 = @(11.53252353746393, 9.262760797427363, 12.89875651657895, 19.00274580682531, 17.81289437971144, 11.27376491217689, 13.97926558861123, 15.93271089461797, 15.70470667162819)
    10 = @(12.45144173811014, 9.653594971327994, 13.154829100926, 20.67494806633232, 17.99809184883239, 11.21630262871601, 14.47423770674849, 15.88012730581339, 15.65027358961744)
    11 = @(12.8458515835829, 9.8253269156488, 13.27472721889945, 21.3917225636224, 18.09153666158632, 11.19256343070284, 14.69638086799627, 15.85805056059112, 15.6306122344677)
    12 = @(12.99180641463258, 9.889444472006881, 13.32057567655743, 21.65686569030329, 18.12821464556412, 11.18392010019807, 14.77998360282139, 15.84995521445153, 15.62390426439373)
    13 = @(12.96052366180211, 9.875676867861525, 13.31068228829838, 21.60004134736742, 18.12025842258988, 11.18576618389591, 14.76200257237707, 15.85168693003911, 15.62531609344242)
    14 = @(12.85792752140969, 9.830620441897995, 13.27849052574499, 21.4136618050453, 18.09452846452476, 11.19184539504207, 14.70326975755112, 15.85737925066018, 15.63004556695957)
    15 = @(12.79464162264143, 9.802901923682194, 13.25882884138165, 21.29868154950795, 18.07893548012699, 11.19561420001544, 14.6672243576108, 15.86090041052485, 15.6330386412914)
    16 = @(12.42519055657699, 9.642246235865166, 13.14705828498123, 20.62722412089977, 17.9921675999007, 11.2179024325339, 14.45965206079619, 15.88160713213789, 15.65166148396547)
    17 = @(12.19249268020163, 9.54210562514705, 13.07933239654836, 20.20408069597325, 17.94127381098462, 11.2321911106746, 14.33147809816333, 15.89478191784884, 15.66439526449047)
    18 = @(12.05642866972088, 9.483939194064854, 13.04070189790573, 19.95656407809801, 17.91286882387438, 11.24063548666158, 14.25747603593115, 15.90253328168449, 15.67219927390366)
    19 = @(12.00997824013075, 9.464148803705916, 13.02767920787671, 19.87204792380568, 17.90340127977939, 11.24353337702059, 14.23237478382185, 15.90518758738594, 15.67492387115487)
    20 = @(12.21749381456287, 9.552824889171346, 13.08650875834178, 20.24955283636154, 17.94660193652818, 11.23064666894127, 14.34515201363729, 15.89336147862653, 15.66299003693604)
    21 = @(12.88815477696565, 9.843879716166736, 13.28793429108671, 21.46857628470577, 18.10205115867953, 11.19005038022201, 14.72053570830185, 15.85570009981669, 15.62863636501601)
    22 = @(13.30665525184325, 10.02875906337189, 13.42215397621811, 22.22866616901552, 18.21116310179751, 11.16553665179177, 14.96281600260023, 15.83262850185006, 15.61048367984159)
    23 = @(13.08510733375986, 9.930587092244881, 13.3502975557878, 21.82633154458858, 18.15225112010189, 11.17843509829595, 14.83381217213777, 15.84480128982254, 15.61977752139601)
    24 = @(12.20619791211374, 9.547980557684136, 13.08326336857906, 20.22900810905287, 17.94419042692851, 11.2313441958447, 14.3389710054087, 15.8940031077956, 15.66362383549741)
    25 = @(11.17323983891486, 9.113762433334577, 12.80776829266617, 18.34778573295695, 17.75302500047678, 11.29704806051805, 13.79488989979752, 15.95371355708131, 15.7292464187332)
}

$cols = @("B", "C", "E", "F", "G", "H", "M", "N", "O")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
